# B6-PowerPoint.pptx edit:
#  1) Re-style the three tables (slides 14-16) from the plain custom
#     "Table_0" style to the built-in "Medium Style 2 - Accent 1" style.
#  2) Swap the presentation's colour theme ("Integral"/Red Violet) for
#     the plain "Office Theme" colours (font scheme / format scheme are
#     identical between the two themes, so only the 12 scheme colours
#     actually change).

$p = $ppt.ActivePresentation

# --- 1) Table styles ---------------------------------------------------
$newTableStyle = "{7995624E-0189-4BD7-ADA7-A8BCA74A5438}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $tableShape = $slide.Shapes.Item(1)
    $tableShape.Table.ApplyStyle($newTableStyle, $false)
}

# --- 2) Theme colours ---------------------------------------------------
# RGB() in VBA/COM packs as 0xBBGGRR, so the literals below are the
# byte-swapped form of the target "Office Theme" hex colours.
$officeThemeColors = @{
    1  = 0x000000  # dk1      000000
    2  = 0xFFFFFF  # lt1      FFFFFF
    3  = 0x6A5444  # dk2      44546A
    4  = 0xE6E6E7  # lt2      E7E6E6
    5  = 0xD59B5B  # accent1  5B9BD5
    6  = 0x317DED  # accent2  ED7D31
    7  = 0xA5A5A5  # accent3  A5A5A5
    8  = 0x00C0FF  # accent4  FFC000
    9  = 0xC47244  # accent5  4472C4
    10 = 0x47AD70  # accent6  70AD47
    11 = 0xC16305  # hlink    0563C1
    12 = 0x724F95  # folHlink 954F72
}

$tcs = $p.Slides.Item(1).ThemeColorScheme
foreach ($i in 1..12) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i]
}
